$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:K data to E:L
$ws.Range("D1").EntireColumn.Insert()

# Copy the formatting (now sitting in column E, the old column D) into the new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clean up stray blank cells created by the insert/paste in rows that never had column-D data
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# Populate the new column D with the latest (FY2018) financial data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 885600
$ws.Range("D9").Value = 774400
$ws.Range("D10").Value = 111200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 600
$ws.Range("D17").Value = 803200
$ws.Range("D18").Value = 82400
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 193500
$ws.Range("D22").Value = 11600
$ws.Range("D23").Value = 70800
$ws.Range("D24").Value = 7000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 63800
$ws.Range("D27").Value = 63800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 63800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 63800
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 172500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 307900
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 1595900
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 219400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 83000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3136800
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 2163300
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 222400
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2427500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 79800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 709200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 63800
$ws.Range("D83").Value = 111100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 290000
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -266800
$ws.Range("D96").Value = -36100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -14300
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 9000
